$d = $word.ActiveDocument

# --- Part 1: remove the _GoBack bookmark from its original location ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Part 2: append a new paragraph at the end of the document body ---
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($n)
$lastRange = $lastPara.Range
$lastRange.Collapse(0)
$lastRange.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newRange = $newPara.Range
$newRange.Collapse(0)

$xmlFragment = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Le dernier point bloquant c’est la découverte de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>heroku</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> et de l’apprentissage de la mise en ligne de l’application. Cela souleva quelques soucis dut à l’utilisation de la b</w:t></w:r><w:r><w:t xml:space="preserve">ibliothèque </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>nltk</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> qui demanda des configurations particulière. Une fois mise en place, l’application fut mise en ligne assez simplement.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$newRange.InsertXML($xmlFragment)
